$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the existing row 401, shifting the current
# rows 401:427 down to 404:430 (dates / prices further down the sheet).
$ws.Rows("401:403").Insert()

# Common (unchanging) column values for this data block.
$colA = 6
$colB = "Mercado Mayorista Lo Valledor de Santiago"
$colC = "Metropolitana"
$colE = 13
$colF = "Fruta"
$colG = 100107
$colH = "Otros"
$colI = 100107002
$colJ = "Chirimoya"
$colK = "Cultivar IV Región"

# New row data: Date, Grade(L), M, N, O, P, Q(unit), R(province), S, T
$newRows = @(
    @{ Row = 401; D = 45166; L = "Especial"; M = 100; N = 28000; O = 28000; P = 28000; Q = "`$/bandeja 10 kilos"; R = "Provincia de Limarí"; S = 2800; T = 10 },
    @{ Row = 402; D = 45166; L = "Primera";  M = 200; N = 26000; O = 26000; P = 26000; Q = "`$/bandeja 10 kilos"; R = "Provincia de Limarí"; S = 2600; T = 10 },
    @{ Row = 403; D = 45166; L = "Segunda";  M = 220; N = 23000; O = 23000; P = 23000; Q = "`$/bandeja 10 kilos"; R = "Provincia de Limarí"; S = 2300; T = 10 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $colA
    $ws.Cells.Item($row, 2).Value = $colB
    $ws.Cells.Item($row, 3).Value = $colC
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $colE
    $ws.Cells.Item($row, 6).Value = $colF
    $ws.Cells.Item($row, 7).Value = $colG
    $ws.Cells.Item($row, 8).Value = $colH
    $ws.Cells.Item($row, 9).Value = $colI
    $ws.Cells.Item($row, 10).Value = $colJ
    $ws.Cells.Item($row, 11).Value = $colK
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
